# RosterValidationErrorTool.xlsx - "protected sheets and minor instruction changes"
#
# 1. Clear the leftover sample/test roster-error rows (B3:F31) on the
#    "Errors" sheet. This was scratch data (RF / ow4270302 / ... ) that
#    leaked into the template; removing it also drops the now-unused
#    shared strings and re-indexes the ones that follow.
# 2. Protect the "Instructions" and "helpText" sheets.
# 3. Restore the various sheet selections / active-cell state and make
#    "Dashboard" the active tab again.

$wb = $excel.ActiveWorkbook

# --- 1. Clear the sample error rows on the Errors sheet -------------------
$errors = $wb.Worksheets.Item("Errors")
$errors.Range("B3:F31").ClearContents()

# --- 2. Protect the Instructions and helpText sheets -----------------------
$instructions = $wb.Worksheets.Item("Instructions")
$instructions.Protect("nwea")

$helpText = $wb.Worksheets.Item("helpText")
$helpText.Protect("nwea")

# --- 3. Selections on each sheet (order matters: last-selected sheet
#        becomes the active tab, so Dashboard is done last) ---------------
$errors.Range("I21").Select()

$calcs = $wb.Worksheets.Item("Calcs")
$calcs.Range("B75").Select()

$notes = $wb.Worksheets.Item("NOTES")
$notes.Range("C34").Select()

$instructions.Range("H30").Select()

$helpText.Range("C12").Select()

$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Activate()
$dashboard.Range("F10").Select()
